# Applies the cryptocurrency price/volume refresh described by the commit diff.
# Updates columns D (Price) and E (Volume(1h)) for most rows, and for rows 41/42
# also swaps the Coin name (B) and Link (C) so RenderToken/Stacks trade places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.562.09'
$ws.Range('E2').Value = '  -2.26%  '
# Row 3
$ws.Range('D3').Value = '2.644.41'
$ws.Range('E3').Value = '  -3.48%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D5').Value = '598.75'
$ws.Range('E5').Value = '  -1.26%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D6').Value = '168.17'
$ws.Range('E6').Value = '  -1.41%  '
# Row 7
$ws.Range('E7').Value = '  +0.03%  '
# Row 8
$ws.Range('E8').Value = '  -0.89%  '
# Row 9
$ws.Range('D9').Value = '2.644.80'
$ws.Range('E9').Value = '  -3.40%  '
# Row 10
$ws.Range('E10').Value = '  -1.74%  '
# Row 11
$ws.Range('E11').Value = '  +1.88%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D12').Value = '0.366'
$ws.Range('E12').Value = '  -1.46%  '
# Row 13
$ws.Range('E13').Value = '  -2.22%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D14').Value = '28.07'
$ws.Range('E14').Value = '  -3.03%  '
# Row 15
$ws.Range('D15').Value = '3.118.86'
$ws.Range('E15').Value = '  -3.67%  '
# Row 16
$ws.Range('E16').Value = '  -4.03%  '
# Row 17
$ws.Range('D17').Value = '67.446.52'
$ws.Range('E17').Value = '  -2.31%  '
# Row 18
$ws.Range('D18').Value = '2.644.58'
$ws.Range('E18').Value = '  -2.86%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D19').Value = '11.89'
$ws.Range('E19').Value = '  -0.61%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D20').Value = '7.90'
$ws.Range('E20').Value = '  +2.44%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D21').Value = '363.77'
$ws.Range('E21').Value = '  -3.13%  '
# Row 22
$ws.Range('E22').Value = '  -3.42%  '
# Row 23
$ws.Range('E23').Value = '  -4.53%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D24').Value = '11.00'
$ws.Range('E24').Value = '  +8.26%  '
# Row 25
$ws.Range('E25').Value = '  -5.52%  '
# Row 26
$ws.Range('E26').Value = '  +0.05%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D27').Value = '70.77'
$ws.Range('E27').Value = '  -4.40%  '
# Row 28
$ws.Range('D28').Value = '2.777.47'
$ws.Range('E28').Value = '  -3.31%  '
# Row 29
$ws.Range('E29').Value = '  -4.23%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  -0.28%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D31').Value = '558.10'
$ws.Range('E31').Value = '  -5.46%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D32').Value = '8.03'
$ws.Range('E32').Value = '  -4.28%  '
# Row 33
$ws.Range('E33').Value = '  -4.42%  '
# Row 34
$ws.Range('E34').Value = '  -2.79%  '
# Row 35
$ws.Range('E35').Value = '  -0.16%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.07%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D37').Value = '1.55'
$ws.Range('E37').Value = '  -5.55%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D38').Value = '157.92'
$ws.Range('E38').Value = '  -2.81%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D39').Value = '19.41'
$ws.Range('E39').Value = '  -3.56%  '
# Row 40
$ws.Range('E40').Value = '  -2.75%  '
# Row 41
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D41').Value = '5.29'
$ws.Range('E41').Value = '  -4.30%  '
# Row 42
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D42').Value = '1.83'
$ws.Range('E42').Value = '  -5.07%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D43').Value = '17.93'
$ws.Range('E43').Value = '  -0.44%  '
# Row 44
$ws.Range('E44').Value = '  -6.05%  '
# Row 45
$ws.Range('E45').Value = '  +0.01%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D46').Value = '40.21'
$ws.Range('E46').Value = '  -2.22%  '
# Row 47
$ws.Range('E47').Value = '  -3.46%  '
# Row 48
$ws.Range('E48').Value = '  -1.92%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'  # keep as text so trailing zeros/decimals survive
$ws.Range('D49').Value = '153.88'
$ws.Range('E49').Value = '  -1.85%  '
# Row 51
$ws.Range('E51').Value = '  -4.36%  '
